# "Security Attribute now executing; de-coupled Registration of critical
# functionalities."
#
# Practically, this commit:
#   1. Adds a new "Essential" TODO item under the "Generic Backlog" sheet:
#        "Goto the Doctor and get checked up + Get my shoulder checked"
#      (the old, shorter "Goto the Doctor and get checked up" row is removed),
#   2. Removes the "Bike Rack, Inner Tube, Pumps" TODO line, and
#   3. Tweaks the wording of the "T-Shirts..." line to
#        "T-Shirts (get them!), Jeans, Kicks, Jack Purcells"
#   4. Leaves "Generic Backlog" as the active/selected sheet (with A10
#      selected) instead of "SCIGON Payroll Summary".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Generic Backlog")

# --- Insert the new "Goto the doctor ..." row above the current row 16 ---
# Row 16 currently holds "House Hunting Backlog"; push it (and everything
# below) down by one and give the fresh row the same yellow "TODO" look
# used elsewhere on this sheet (e.g. row 8 "De-commission Big Blue").
$ws.Rows("16:16").Insert()
$ws.Range("A8:B8").Copy()
$ws.Range("A16:B16").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A16").Value = "Goto the Doctor and get checked up + Get my shoulder checked"
$ws.Range("B16").Value = "TODO"

# The original (shorter) "Goto the Doctor and get checked up" row has now
# been pushed down to row 18 - remove it outright, it's superseded by the
# row we just inserted above.
$ws.Rows("18:18").Delete()

# --- Drop the "Bike Rack, Inner Tube, Pumps" line entirely ---
$ws.Rows("21:21").Delete()

# --- Reword the T-Shirts line (now sitting at row 22) ---
$ws.Range("A22").Value = "T-Shirts (get them!), Jeans, Kicks, Jack Purcells"

# --- Make "Generic Backlog" the active sheet/selection on save ---
$ws.Activate()
$ws.Range("A10").Select() | Out-Null
